$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H12").Value = 1000
$ws_ALC.Range("I12").Value = 1000
$ws_ALC.Range("J12").Value = 0
$ws_ALC.Range("K12").Value = 1000
$ws_ALC.Range("L12").Value = 0
$ws_ALC.Range("M12").Value = -830
$ws_ALC.Range("N12").ClearContents()
$ws_ALC.Range("J21").Value = 27500
$ws_ALC.Range("L21").Value = 27500
$ws_ALC.Range("N21").Value = -28436
$ws_ALC.Range("J23").Value = 27500
$ws_ALC.Range("L23").Value = 27500
$ws_ALC.Range("N23").Value = -27968
$ws_ALC.Range("H58").Value = 2298.8462
$ws_ALC.Range("I58").Value = 228.75
$ws_ALC.Range("J58").Value = 3218.889
$ws_ALC.Range("K58").Value = 686.25
$ws_ALC.Range("L58").Value = 9656.667000000001
$ws_ALC.Range("M58").Value = -536.25
$ws_ALC.Range("N58").Value = -9956.667000000001
$ws_ALC.Range("H87").Value = 100077
$ws_ALC.Range("J87").Value = 100077
$ws_ALC.Range("L87").Value = 100077
$ws_ALC.Range("N87").Value = -102573
$ws_ALC.Range("H90").Value = 100077
$ws_ALC.Range("J90").Value = 100077
$ws_ALC.Range("L90").Value = 300231
$ws_ALC.Range("N90").Value = -312711
$ws_ALC.Range("H112").Value = 2240.1365
$ws_ALC.Range("J112").Value = 2369.125
$ws_ALC.Range("L112").Value = 7107.375
$ws_ALC.Range("N112").Value = -9323.375
$ws_ALC.Range("H129").Value = 955.7646999999999
$ws_ALC.Range("I129").Value = 480.26666
$ws_ALC.Range("J129").Value = 1153.8889
$ws_ALC.Range("K129").Value = 1440.79998
$ws_ALC.Range("L129").Value = 3461.6667
$ws_ALC.Range("M129").Value = 3559.20002
$ws_ALC.Range("N129").Value = -13461.6667
$ws_ALC.Range("H133").Value = 37490
$ws_ALC.Range("J133").Value = 37490
$ws_ALC.Range("L133").Value = 37490
$ws_ALC.Range("N133").Value = -47610
$ws_ALC.Range("H136").Value = 59857.777
$ws_ALC.Range("J136").Value = 59857.777
$ws_ALC.Range("L136").Value = 59857.777
$ws_ALC.Range("N136").Value = -70057.777
$ws_ALC.Range("H138").Value = 3453522.8
$ws_ALC.Range("I138").Value = 9527031
$ws_ALC.Range("J138").Value = 6396.3516
$ws_ALC.Range("K138").Value = 28581093
$ws_ALC.Range("L138").Value = 19189.0548
$ws_ALC.Range("M138").Value = -28575953
$ws_ALC.Range("N138").Value = -29469.0548
$ws_ALC.Range("H139").Value = 80000
$ws_ALC.Range("J139").Value = 80000
$ws_ALC.Range("L139").Value = 80000
$ws_ALC.Range("N139").Value = -90280
$ws_ALC.Range("H140").Value = 110260
$ws_ALC.Range("J140").Value = 113602
$ws_ALC.Range("L140").Value = 113602
$ws_ALC.Range("N140").Value = -123962
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H138").Value = 61471.668
$ws_ARM.Range("J138").Value = 61471.668
$ws_ARM.Range("L138").Value = 61471.668
$ws_ARM.Range("N138").Value = -71751.66800000001
$ws_ARM.Range("H139").Value = 89177.125
$ws_ARM.Range("J139").Value = 89177.125
$ws_ARM.Range("L139").Value = 89177.125
$ws_ARM.Range("N139").Value = -99457.125
$ws_ARM.Range("H140").Value = 49533.207
$ws_ARM.Range("J140").Value = 49533.207
$ws_ARM.Range("L140").Value = 49533.207
$ws_ARM.Range("N140").Value = -59893.207
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H105").Value = 3715.4
$ws_BSM.Range("I105").Value = 3266.5
$ws_BSM.Range("K105").Value = 3266.5
$ws_BSM.Range("M105").Value = -1519.5
$ws_BSM.Range("H132").Value = 76313.08
$ws_BSM.Range("J132").Value = 76313.08
$ws_BSM.Range("L132").Value = 76313.08
$ws_BSM.Range("N132").Value = -86433.08
$ws_BSM.Range("H137").Value = 45567.273
$ws_BSM.Range("J137").Value = 45567.273
$ws_BSM.Range("L137").Value = 45567.273
$ws_BSM.Range("N137").Value = -55767.273
$ws_BSM.Range("H138").Value = 51883.332
$ws_BSM.Range("J138").Value = 51883.332
$ws_BSM.Range("L138").Value = 51883.332
$ws_BSM.Range("N138").Value = -62163.332
$ws_BSM.Range("H140").Value = 49309.41
$ws_BSM.Range("J140").Value = 49309.41
$ws_BSM.Range("L140").Value = 49309.41
$ws_BSM.Range("N140").Value = -59669.41
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H99").Value = 2928.3872
$ws_CRP.Range("I99").Value = 2876.4707
$ws_CRP.Range("J99").Value = 2991.4285
$ws_CRP.Range("K99").Value = 2876.4707
$ws_CRP.Range("L99").Value = 2991.4285
$ws_CRP.Range("M99").Value = -1378.4707
$ws_CRP.Range("N99").Value = -5987.4285
$ws_CRP.Range("H126").Value = 2928.3872
$ws_CRP.Range("I126").Value = 2876.4707
$ws_CRP.Range("J126").Value = 2991.4285
$ws_CRP.Range("K126").Value = 8629.4121
$ws_CRP.Range("L126").Value = 8974.2855
$ws_CRP.Range("M126").Value = -6159.4121
$ws_CRP.Range("N126").Value = -13914.2855
$ws_CRP.Range("H135").Value = 87692
$ws_CRP.Range("J135").Value = 121538
$ws_CRP.Range("L135").Value = 121538
$ws_CRP.Range("N135").Value = -131678
$ws_CRP.Range("H137").Value = 54233.332
$ws_CRP.Range("J137").Value = 73620
$ws_CRP.Range("L137").Value = 73620
$ws_CRP.Range("N137").Value = -83820
$ws_CRP.Range("H138").Value = 58816.43
$ws_CRP.Range("J138").Value = 58816.43
$ws_CRP.Range("L138").Value = 58816.43
$ws_CRP.Range("N138").Value = -69096.42999999999
$ws_CRP.Range("H140").Value = 72828.75
$ws_CRP.Range("J140").Value = 72828.75
$ws_CRP.Range("L140").Value = 72828.75
$ws_CRP.Range("N140").Value = -83188.75
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H129").Value = 2175781.2
$ws_CUL.Range("I129").Value = 825
$ws_CUL.Range("J129").Value = 2633666.8
$ws_CUL.Range("K129").Value = 2475
$ws_CUL.Range("L129").Value = 7901000.399999999
$ws_CUL.Range("M129").Value = 2525
$ws_CUL.Range("N129").Value = -7911000.399999999
$ws_CUL.Range("H131").Value = 13160708
$ws_CUL.Range("J131").Value = 15153037
$ws_CUL.Range("L131").Value = 45459111
$ws_CUL.Range("N131").Value = -45469191
$ws_CUL.Range("H140").Value = 2611.037
$ws_CUL.Range("I140").Value = 1073.1428
$ws_CUL.Range("K140").Value = 3219.4284
$ws_CUL.Range("M140").Value = 1960.5716
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H133").Value = 62780
$ws_GSM.Range("J133").Value = 62780
$ws_GSM.Range("L133").Value = 62780
$ws_GSM.Range("N133").Value = -72900
$ws_GSM.Range("H138").Value = 59770
$ws_GSM.Range("J138").Value = 59770
$ws_GSM.Range("L138").Value = 59770
$ws_GSM.Range("N138").Value = -70050
$ws_GSM.Range("H140").Value = 50654
$ws_GSM.Range("J140").Value = 50654
$ws_GSM.Range("L140").Value = 50654
$ws_GSM.Range("N140").Value = -61014
$ws_GSM.Range("H141").Value = 56173.75
$ws_GSM.Range("J141").Value = 56173.75
$ws_GSM.Range("L141").Value = 56173.75
$ws_GSM.Range("N141").Value = -66533.75
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 468.75
$ws_LTW.Range("I22").Value = 507.14285
$ws_LTW.Range("K22").Value = 507.14285
$ws_LTW.Range("M22").Value = -212.14285
$ws_LTW.Range("H27").Value = 468.75
$ws_LTW.Range("I27").Value = 507.14285
$ws_LTW.Range("K27").Value = 507.14285
$ws_LTW.Range("M27").Value = -400.14285
$ws_LTW.Range("H40").Value = 5782.7085
$ws_LTW.Range("I40").Value = 5730
$ws_LTW.Range("J40").Value = 5888.125
$ws_LTW.Range("K40").Value = 5730
$ws_LTW.Range("L40").Value = 5888.125
$ws_LTW.Range("M40").Value = -5594
$ws_LTW.Range("N40").Value = -6160.125
$ws_LTW.Range("H46").Value = 1667.1111
$ws_LTW.Range("I46").Value = 1500
$ws_LTW.Range("J46").Value = 1688
$ws_LTW.Range("K46").Value = 1500
$ws_LTW.Range("L46").Value = 1688
$ws_LTW.Range("M46").Value = -1312
$ws_LTW.Range("N46").Value = -2064
$ws_LTW.Range("H133").Value = 55625.04
$ws_LTW.Range("J133").Value = 55625.04
$ws_LTW.Range("L133").Value = 55625.04
$ws_LTW.Range("N133").Value = -60685.04
$ws_LTW.Range("H135").Value = 163388.17
$ws_LTW.Range("J135").Value = 163388.17
$ws_LTW.Range("L135").Value = 163388.17
$ws_LTW.Range("N135").Value = -173528.17
$ws_LTW.Range("H137").Value = 51199.715
$ws_LTW.Range("J137").Value = 51199.715
$ws_LTW.Range("L137").Value = 51199.715
$ws_LTW.Range("N137").Value = -61399.715
$ws_LTW.Range("H139").Value = 79895
$ws_LTW.Range("J139").Value = 79895
$ws_LTW.Range("L139").Value = 79895
$ws_LTW.Range("N139").Value = -90175
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H125").Value = 65635.55499999999
$ws_WVR.Range("J125").Value = 65635.55499999999
$ws_WVR.Range("L125").Value = 65635.55499999999
$ws_WVR.Range("N125").Value = -75475.55499999999
$ws_WVR.Range("H135").Value = 59313.89
$ws_WVR.Range("J135").Value = 59313.89
$ws_WVR.Range("L135").Value = 59313.89
$ws_WVR.Range("N135").Value = -69453.89
$ws_WVR.Range("H137").Value = 48000
$ws_WVR.Range("J137").Value = 48000
$ws_WVR.Range("L137").Value = 48000
$ws_WVR.Range("N137").Value = -58200
$ws_WVR.Range("H139").Value = 60587.5
$ws_WVR.Range("J139").Value = 60587.5
$ws_WVR.Range("L139").Value = 60587.5
$ws_WVR.Range("N139").Value = -70867.5
$ws_WVR.Range("H140").Value = 43220
$ws_WVR.Range("J140").Value = 43220
$ws_WVR.Range("L140").Value = 43220
$ws_WVR.Range("N140").Value = -53580
